$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 4 (Nov20th-Nov27th) block: new task row, then the Weekly Total row moves up ---
$ws.Range("B28").Value = "Looking into MIME/HTTP Post"
$ws.Range("C28").Value = 0.0625
$ws.Range("C28").NumberFormat = "h:mm"

$ws.Range("B29").Value = "Weekly Total"
$ws.Range("C29").Formula = "=SUM(C25:C28)"
$ws.Range("C29").NumberFormat = "h:mm"

# --- New week block: Week:Nov27th-Dec3rd (rows 30-35) ---
$ws.Range("A30").Value = "Week:Nov27th-Dec3rd"

$ws.Range("B31").Value = "Task"
$ws.Range("C31").Clear()

$ws.Range("B32").Value = "Meetings"
$ws.Range("C32").Value = 0.08333333333333333
$ws.Range("C32").NumberFormat = "h:mm"

$ws.Range("B33").Value = "Write up of algorithms"
$ws.Range("C33").Value = 0.041666666666666664
$ws.Range("C33").NumberFormat = "h:mm"

$ws.Range("B34").Value = "Research into algorithms"
$ws.Range("C34").Value = 0.10416666666666667
$ws.Range("C34").NumberFormat = "h:mm"

$ws.Range("B35").Value = "Weekly Total"
$ws.Range("C35").Formula = "=SUM(C32:C34)"
$ws.Range("C35").NumberFormat = "h:mm"

# --- Blank formatted filler rows (time-format style carried down, matching the sheet's style) ---
$ws.Range("C36").NumberFormat = "h:mm"
$ws.Range("C37").NumberFormat = "h:mm"
$ws.Range("C38").NumberFormat = "h:mm"
$ws.Range("C39").NumberFormat = "h:mm"
$ws.Range("C40").NumberFormat = "h:mm"
$ws.Range("D40").NumberFormat = "h:mm"
$ws.Range("C41").NumberFormat = "h:mm"
$ws.Range("C42").NumberFormat = "h:mm"
$ws.Range("C43").NumberFormat = "h:mm"

# --- Current Total moves from row 34 -> row 44; formula now sums the 5 weekly totals with commas ---
$ws.Range("B44").Value = "Current Total:"
$ws.Range("C44").Formula = "=SUM(C9,C16,C22,C29,C35)"
$ws.Range("C44").NumberFormat = "h:mm"

# --- View / selection state ---
$ws.Range("D41").Select()
